# Generate Report for Handoff
# Adds two new localization entries -
#   255e04c5-ced2-4fbc-8285-a86bd0231d73  (inserted before the existing 5261da22 row)
#   66ec3868-ad0c-45fa-9a94-6068a8de1ec8  (appended after the existing 5261da22 row)
# to the Overview / zh-cn / de-de sheets, pushing 5261da22's row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop every existing hyperlink up front - cell values get rewritten below and
# stale hyperlink bindings do not follow a plain .Value re-assignment.
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = "139671fe-2709-4645-a438-38e543dc5459.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-32-18 20:32:19"

$ws1.Range("A3").Value = "255e04c5-ced2-4fbc-8285-a86bd0231d73.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-33-18 20:33:10"

$ws1.Range("A4").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-32-18 20:32:01"

$ws1.Range("A5").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-33-18 20:33:10"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3106469f70392444d58b2e7f959cfdf18a1a2a98/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/e2e/255e04c5-ced2-4fbc-8285-a86bd0231d73.md", "", "", "255e04c5-ced2-4fbc-8285-a86bd0231d73.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/868691f7373919591d3a10d4348b4462bbbb79c7/e2e/5261da22-23e2-4c5f-a60b-446bc987709a.md", "", "", "5261da22-23e2-4c5f-a60b-446bc987709a.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/e2e/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md", "", "", "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = "139671fe-2709-4645-a438-38e543dc5459.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-18 20:32:16"
$ws2.Range("F2").Value = "139671fe-2709-4645-a438-38e543dc5459.md"
$ws2.Range("G2").Value = "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-18 20:32:34"
$ws2.Range("I2").Value = "Include"

$ws2.Range("A3").Value = "255e04c5-ced2-4fbc-8285-a86bd0231d73.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-18 20:33:07"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$ws2.Range("A4").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-18 20:31:58"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "Include"

$ws2.Range("A5").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-18 20:33:07"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("I5").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3106469f70392444d58b2e7f959cfdf18a1a2a98/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3106469f70392444d58b2e7f959cfdf18a1a2a98/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03c7cc47042e45c1722ca135221ff641cfad956d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf", "", "", "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3bd8d5bf555c6e26790de318d6d30cf6add3b64c/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3a9f76515a9d03f33b4b29db357c3927dc2d4e8f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf", "", "", "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/e2e/255e04c5-ced2-4fbc-8285-a86bd0231d73.md", "", "", "255e04c5-ced2-4fbc-8285-a86bd0231d73.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/e2e/255e04c5-ced2-4fbc-8285-a86bd0231d73.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.zh-cn.xlf", "", "", "255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/868691f7373919591d3a10d4348b4462bbbb79c7/e2e/5261da22-23e2-4c5f-a60b-446bc987709a.md", "", "", "5261da22-23e2-4c5f-a60b-446bc987709a.md")
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/868691f7373919591d3a10d4348b4462bbbb79c7/e2e/5261da22-23e2-4c5f-a60b-446bc987709a.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2c7a51b5e515d441135facd375e43fa37cf70d0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.zh-cn.xlf", "", "", "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/e2e/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md", "", "", "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md")
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/e2e/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.zh-cn.xlf", "", "", "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = "139671fe-2709-4645-a438-38e543dc5459.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-18 20:32:19"
$ws3.Range("F2").Value = "139671fe-2709-4645-a438-38e543dc5459.md"
$ws3.Range("G2").Value = "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-18 20:32:39"
$ws3.Range("I2").Value = "Include"

$ws3.Range("A3").Value = "255e04c5-ced2-4fbc-8285-a86bd0231d73.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-18 20:33:10"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$ws3.Range("A4").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-18 20:32:01"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "Include"

$ws3.Range("A5").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-18 20:33:10"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("I5").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3106469f70392444d58b2e7f959cfdf18a1a2a98/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3106469f70392444d58b2e7f959cfdf18a1a2a98/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce6a700e80eb2c6ec27caa4d6c786a4b09b6e930/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf", "", "", "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ec0e86947ef9c6f9de1ec88c2b56029e5136be08/e2e/139671fe-2709-4645-a438-38e543dc5459.md", "", "", "139671fe-2709-4645-a438-38e543dc5459.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27af7ec660eb2e64ef4536410940ffa6c989a655/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf", "", "", "139671fe-2709-4645-a438-38e543dc5459.e762053e4d4d7c224bbc13ac0accb9938df9dbbf.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/e2e/255e04c5-ced2-4fbc-8285-a86bd0231d73.md", "", "", "255e04c5-ced2-4fbc-8285-a86bd0231d73.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/e2e/255e04c5-ced2-4fbc-8285-a86bd0231d73.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410d2d1b10eed83d873f2e4e166a65f7ce3b074/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.de-de.xlf", "", "", "255e04c5-ced2-4fbc-8285-a86bd0231d73.8410d2d1b10eed83d873f2e4e166a65f7ce3b074.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/868691f7373919591d3a10d4348b4462bbbb79c7/e2e/5261da22-23e2-4c5f-a60b-446bc987709a.md", "", "", "5261da22-23e2-4c5f-a60b-446bc987709a.md")
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/868691f7373919591d3a10d4348b4462bbbb79c7/e2e/5261da22-23e2-4c5f-a60b-446bc987709a.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9f5d0058f316f12445ad398008d233dd6705f2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.de-de.xlf", "", "", "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/e2e/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md", "", "", "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md")
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/e2e/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b926bccb71c10874a9e9887bbb707a64496b2e33/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.de-de.xlf", "", "", "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.de-de.xlf")
